$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0.1812636768995252
$ws.Range("B2").Value = 0.02497434187251466
$ws.Range("E2").Value = 0.1738488
$ws.Range("G2").Value = 0.130401110924859
$ws.Range("H2").Value = 0.4076228425764395
$ws.Range("I2").Value = 0.9948565999999999
$ws.Range("N2").Value = 8.552054364790164
$ws.Range("O2").Value = 5.606994352318015

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.2169040844903429
$ws.Range("B2").Value = 0.1832692253478948
$ws.Range("E2").Value = 0.1512647663897953
$ws.Range("H2").Value = 0.0886774357781539
$ws.Range("I2").Value = 0.8712850990390117
$ws.Range("N2").Value = 5.772063423759102
$ws.Range("O2").Value = 4.071807768933123

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.3694649269708148
$ws.Range("B2").Value = 0.000003555464936272346
$ws.Range("I2").Value = 0.4843915289694924
$ws.Range("M2").Value = 0.02289547199149599
$ws.Range("N2").Value = 4.990247716527666
$ws.Range("O2").Value = 6.450188393386924
